$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 35
$ws.Range("F6").Value = 25.85023498535156

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 28
$ws.Range("D7").Value = 41.47622680664062
